$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create TestData02 (API Positive) ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "TestData02"
$ws1.Range("A1:F2").Copy($ws2.Range("A1"))
for ($i = 1; $i -le 6; $i++) {
    $ws2.Columns.Item($i).ColumnWidth = $ws1.Columns.Item($i).ColumnWidth
}
$ws2.Range("A2").Value = "API Positive"
$ws2.Range("E2").Value = "eve.holt@reqres.in"
$ws2.Hyperlinks.Add($ws2.Range("E2"), "mailto:eve.holt@reqres.in")
$ws2.Range("F2").Value = "cityslicka"
$ws2.PageSetup.Orientation = 1
$ws2.Range("F2").Select() | Out-Null

# --- Create TestData03 (API Negative) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "TestData03"
$ws1.Range("A1:F2").Copy($ws3.Range("A1"))
for ($i = 1; $i -le 6; $i++) {
    $ws3.Columns.Item($i).ColumnWidth = $ws1.Columns.Item($i).ColumnWidth
}
$ws3.Range("A2").Value = "API Negative"
$ws3.Range("E2").Value = "test@gmail.com"
$ws3.Hyperlinks.Add($ws3.Range("E2"), "mailto:test@gmail.com")
$ws3.Range("F2").Value = "test"
$ws3.PageSetup.Orientation = 1
$ws3.Range("E3").Select() | Out-Null
